$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.678.89"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").Value = "'1.744.72"
$ws.Range("E3").Value = "  -5.70%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'236.30"
$ws.Range("E5").Value = "  -9.93%  "
$ws.Range("E7").Value = "  -8.78%  "
$ws.Range("D8").Value = "'41.53"
$ws.Range("E8").Value = "  -8.00%  "
$ws.Range("D9").Value = "'0.2484"
$ws.Range("E9").Value = "  -22.25%  "
$ws.Range("D10").Value = "'0.05969"
$ws.Range("E10").Value = "  -15.41%  "
$ws.Range("D11").Value = "'1.750.93"
$ws.Range("E11").Value = "  -5.39%  "
$ws.Range("D12").Value = "'0.06798"
$ws.Range("E12").Value = "  -13.19%  "
$ws.Range("D14").Value = "'4.453"
$ws.Range("E14").Value = "  -11.95%  "
$ws.Range("D15").Value = "'77.08"
$ws.Range("E15").Value = "  -13.95%  "
$ws.Range("D16").Value = "'0.5587"
$ws.Range("E16").Value = "  -28.07%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'25.719.21"
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("D20").Value = "'11.39"
$ws.Range("E20").Value = "  -19.65%  "
$ws.Range("D21").Value = "'0.000006514"
$ws.Range("E21").Value = "  -18.80%  "
$ws.Range("D22").Value = "'1.965.90"
$ws.Range("E22").Value = "  -5.93%  "
$ws.Range("D23").Value = "'3.978"
$ws.Range("E23").Value = "  -14.40%  "
$ws.Range("D24").Value = "'5.007"
$ws.Range("E24").Value = "  -17.25%  "
$ws.Range("D25").Value = "'7.835"
$ws.Range("E25").Value = "  -16.82%  "
$ws.Range("D26").Value = "'136.15"
$ws.Range("E26").Value = "  -4.76%  "
$ws.Range("D27").Value = "'1.486"
$ws.Range("E27").Value = "  -12.44%  "
$ws.Range("D28").Value = "'1.801"
$ws.Range("E28").Value = "  -18.89%  "
$ws.Range("D29").Value = "'14.58"
$ws.Range("E29").Value = "  -15.04%  "
$ws.Range("D30").Value = "'101.63"
$ws.Range("D31").Value = "'3.740"
$ws.Range("E31").Value = "  -13.13%  "
$ws.Range("D32").Value = "'0.08008"
$ws.Range("E32").Value = "  -8.59%  "
$ws.Range("D33").Value = "'3.294"
$ws.Range("E33").Value = "  -19.99%  "
$ws.Range("D34").Value = "'0.04382"
$ws.Range("E34").Value = "  -10.27%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'2.605"
$ws.Range("E36").Value = "  -10.13%  "
$ws.Range("D37").Value = "'0.9785"
$ws.Range("E37").Value = "  -14.56%  "
$ws.Range("D38").Value = "'0.6042"
$ws.Range("E38").Value = "  -18.15%  "
$ws.Range("E39").Value = "  -13.95%  "
$ws.Range("D40").Value = "'1.981"
$ws.Range("E40").Value = "  -15.82%  "
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'103.14"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("D43").Value = "'0.01495"
$ws.Range("E43").Value = "  -14.72%  "
$ws.Range("D44").Value = "'0.7544"
$ws.Range("D45").Value = "'5.138"
$ws.Range("E45").Value = "  -13.28%  "
$ws.Range("D46").Value = "'0.3686"
$ws.Range("E46").Value = "  -23.94%  "
$ws.Range("D47").Value = "'0.05199"
$ws.Range("E47").Value = "  -10.95%  "
$ws.Range("D48").Value = "'0.1067"
$ws.Range("E48").Value = "  -14.77%  "
$ws.Range("D49").Value = "'30.01"
$ws.Range("E49").Value = "  -14.60%  "
$ws.Range("D50").Value = "'52.30"
$ws.Range("E50").Value = "  -13.67%  "
$ws.Range("D51").Value = "'5.829"
$ws.Range("E51").Value = "  -24.86%  "